$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 240, shifting existing rows 240:262 down to 241:263.
$ws.Rows.Item(240).Insert()

# Populate the newly inserted row 240 with its data (matching the same
# "Poroto verde / Sin especificar / Primera" series as its neighbours).
$ws.Range("A240").Value = 5
$ws.Range("B240").Value = "Macroferia Regional de Talca"
$ws.Range("C240").Value = "Maule"
$ws.Range("D240").Value = 45106
$ws.Range("E240").Value = 7
$ws.Range("F240").Value = 100112031
$ws.Range("G240").Value = "Poroto verde"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 150
$ws.Range("K240").Value = 23000
$ws.Range("L240").Value = 23000
$ws.Range("M240").Value = 23000
$ws.Range("N240").Value = "$/malla 25 kilos"
$ws.Range("O240").Value = "Perú"
$ws.Range("P240").Value = 920
$ws.Range("Q240").Value = 25
$ws.Range("R240").Value = "Hortaliza"
